$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.237.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.381.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.13%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.31%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -1.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.382.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("E10").Value = "  +2.86%  "

$ws.Range("E11").Value = "  +1.37%  "

$ws.Range("E12").Value = "  +1.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.350"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.07%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.793.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.85%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000167"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.167.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.390.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.52%  "

$ws.Range("E20").Value = "  +1.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.72%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.65%  "

$ws.Range("E25").Value = "  -11.12%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.488.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "508.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0891"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.151"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.33%  "

$ws.Range("E33").Value = "  -3.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.26%  "

$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.57%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.06%  "

$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.380"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "146.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "149.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.32%  "

$ws.Range("E46").Value = "  -0.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.93%  "

$ws.Range("E48").Value = "  +0.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.67%  "

$ws.Range("E50").Value = "  +0.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0911"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.65%  "
